$d = $word.ActiveDocument

$d.Content.Find.Execute("651×9=", $true, $false, $false, $false, $false, $true, 1, $false, "971×2=", 2)
$d.Content.Find.Execute("839×2=", $true, $false, $false, $false, $false, $true, 1, $false, "569×9=", 2)
$d.Content.Find.Execute("232×3=", $true, $false, $false, $false, $false, $true, 1, $false, "732×4=", 2)
$d.Content.Find.Execute("823×8=", $true, $false, $false, $false, $false, $true, 1, $false, "851×9=", 2)
$d.Content.Find.Execute("852×2=", $true, $false, $false, $false, $false, $true, 1, $false, "830×8=", 2)
$d.Content.Find.Execute("521×2=", $true, $false, $false, $false, $false, $true, 1, $false, "936×8=", 2)
$d.Content.Find.Execute("285×5=", $true, $false, $false, $false, $false, $true, 1, $false, "345×6=", 2)
$d.Content.Find.Execute("957×5=", $true, $false, $false, $false, $false, $true, 1, $false, "195×8=", 2)
$d.Content.Find.Execute("934×2=", $true, $false, $false, $false, $false, $true, 1, $false, "569×6=", 2)
$d.Content.Find.Execute("502×6=", $true, $false, $false, $false, $false, $true, 1, $false, "390×9=", 2)
$d.Content.Find.Execute("244×2=", $true, $false, $false, $false, $false, $true, 1, $false, "918×2=", 2)
$d.Content.Find.Execute("368×3=", $true, $false, $false, $false, $false, $true, 1, $false, "686×6=", 2)
$d.Content.Find.Execute("221×3=", $true, $false, $false, $false, $false, $true, 1, $false, "383×4=", 2)
$d.Content.Find.Execute("461×4=", $true, $false, $false, $false, $false, $true, 1, $false, "562×6=", 2)
$d.Content.Find.Execute("443×9=", $true, $false, $false, $false, $false, $true, 1, $false, "247×3=", 2)
$d.Content.Find.Execute("678×4=", $true, $false, $false, $false, $false, $true, 1, $false, "334×3=", 2)
$d.Content.Find.Execute("618×7=", $true, $false, $false, $false, $false, $true, 1, $false, "459×8=", 2)
$d.Content.Find.Execute("258×7=", $true, $false, $false, $false, $false, $true, 1, $false, "898×4=", 2)
$d.Content.Find.Execute("173×9=", $true, $false, $false, $false, $false, $true, 1, $false, "773×3=", 2)
$d.Content.Find.Execute("157×7=", $true, $false, $false, $false, $false, $true, 1, $false, "246×6=", 2)
$d.Content.Find.Execute("465×9=", $true, $false, $false, $false, $false, $true, 1, $false, "477×8=", 2)
$d.Content.Find.Execute("967×8=", $true, $false, $false, $false, $false, $true, 1, $false, "830×4=", 2)
$d.Content.Find.Execute("690×8=", $true, $false, $false, $false, $false, $true, 1, $false, "909×6=", 2)
$d.Content.Find.Execute("947×5=", $true, $false, $false, $false, $false, $true, 1, $false, "267×4=", 2)
$d.Content.Find.Execute("182×8=", $true, $false, $false, $false, $false, $true, 1, $false, "730×2=", 2)
